$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2116.6667
$ws.Range("I40").Value = 2062.5
$ws.Range("J40").Value = 2225
$ws.Range("K40").Value = 2062.5
$ws.Range("L40").Value = 2225
$ws.Range("M40").Value = -1887.5
$ws.Range("N40").Value = -2575
$ws.Range("H64").Value = 6949627
$ws.Range("I64").Value = 15628685
$ws.Range("J64").Value = 6380.3
$ws.Range("K64").Value = 15628685
$ws.Range("L64").Value = 6380.3
$ws.Range("M64").Value = -15628437
$ws.Range("N64").Value = -6876.3
$ws.Range("H67").Value = 6949627
$ws.Range("I67").Value = 15628685
$ws.Range("J67").Value = 6380.3
$ws.Range("K67").Value = 15628685
$ws.Range("L67").Value = 6380.3
$ws.Range("M67").Value = -15627827
$ws.Range("N67").Value = -8096.3
$ws.Range("H74").Value = 4006.6924
$ws.Range("I74").Value = 3576.6
$ws.Range("J74").Value = 4275.5
$ws.Range("K74").Value = 3576.6
$ws.Range("L74").Value = 4275.5
$ws.Range("M74").Value = -2640.6
$ws.Range("N74").Value = -6147.5
$ws.Range("H77").Value = 4006.6924
$ws.Range("I77").Value = 3576.6
$ws.Range("J77").Value = 4275.5
$ws.Range("K77").Value = 17883
$ws.Range("L77").Value = 21377.5
$ws.Range("M77").Value = -13203
$ws.Range("N77").Value = -30737.5
$ws.Range("H141").Value = 2145.9375
$ws.Range("I141").Value = 2145.9375
$ws.Range("K141").Value = 6437.8125
$ws.Range("M141").Value = -1257.8125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 253383.75
$ws.Range("I2").Value = 337374
$ws.Range("J2").Value = 1413
$ws.Range("K2").Value = 337374
$ws.Range("L2").Value = 1413
$ws.Range("M2").Value = -337261
$ws.Range("N2").Value = -1639
$ws.Range("H74").Value = 4164.048
$ws.Range("I74").Value = 849.25
$ws.Range("J74").Value = 14771.4
$ws.Range("K74").Value = 849.25
$ws.Range("L74").Value = 14771.4
$ws.Range("M74").Value = 24.75
$ws.Range("N74").Value = -16519.4
$ws.Range("H77").Value = 4164.048
$ws.Range("I77").Value = 849.25
$ws.Range("J77").Value = 14771.4
$ws.Range("K77").Value = 4246.25
$ws.Range("L77").Value = 73857
$ws.Range("M77").Value = 121.75
$ws.Range("N77").Value = -82593
$ws.Range("H116").Value = 253383.75
$ws.Range("I116").Value = 337374
$ws.Range("J116").Value = 1413
$ws.Range("K116").Value = 337374
$ws.Range("L116").Value = 1413
$ws.Range("M116").Value = -335080
$ws.Range("N116").Value = -6001
$ws.Range("H122").Value = 2663.182
$ws.Range("I122").Value = 2629.5
$ws.Range("K122").Value = 7888.5
$ws.Range("M122").Value = -5438.5
$ws.Range("H127").Value = 40000
$ws.Range("J127").Value = 40000
$ws.Range("L127").Value = 40000
$ws.Range("N127").Value = -49920

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 253383.75
$ws.Range("I3").Value = 337374
$ws.Range("J3").Value = 1413
$ws.Range("K3").Value = 337374
$ws.Range("L3").Value = 1413
$ws.Range("M3").Value = -337260
$ws.Range("N3").Value = -1641
$ws.Range("H86").Value = 1845.091
$ws.Range("I86").Value = 1964.8334
$ws.Range("J86").Value = 1701.4
$ws.Range("K86").Value = 1964.8334
$ws.Range("L86").Value = 1701.4
$ws.Range("M86").Value = -841.8334
$ws.Range("N86").Value = -3947.4
$ws.Range("H89").Value = 1845.091
$ws.Range("I89").Value = 1964.8334
$ws.Range("J89").Value = 1701.4
$ws.Range("K89").Value = 9824.166999999999
$ws.Range("L89").Value = 8507
$ws.Range("M89").Value = -4208.166999999999
$ws.Range("N89").Value = -19739
$ws.Range("H127").Value = 33000
$ws.Range("J127").Value = 33000
$ws.Range("L127").Value = 33000
$ws.Range("N127").Value = -42920

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 16717.941
$ws.Range("I62").Value = 20725.416
$ws.Range("J62").Value = 7100
$ws.Range("K62").Value = 20725.416
$ws.Range("L62").Value = 7100
$ws.Range("M62").Value = -20101.416
$ws.Range("N62").Value = -8348
$ws.Range("H65").Value = 16717.941
$ws.Range("I65").Value = 20725.416
$ws.Range("J65").Value = 7100
$ws.Range("K65").Value = 103627.08
$ws.Range("L65").Value = 35500
$ws.Range("M65").Value = -100507.08
$ws.Range("N65").Value = -41740
$ws.Range("H102").Value = 29750
$ws.Range("J102").Value = 29750
$ws.Range("L102").Value = 29750
$ws.Range("N102").Value = -34618
$ws.Range("H122").Value = 1637.0834
$ws.Range("I122").Value = 1557.3636
$ws.Range("K122").Value = 4672.0908
$ws.Range("M122").Value = -2222.0908

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 565.8
$ws.Range("J107").Value = 545.7143
$ws.Range("L107").Value = 1637.1429
$ws.Range("N107").Value = -5477.1429
$ws.Range("H131").Value = 6804284
$ws.Range("I131").Value = 546
$ws.Range("J131").Value = 7577436
$ws.Range("K131").Value = 1638
$ws.Range("L131").Value = 22732308
$ws.Range("M131").Value = 3402
$ws.Range("N131").Value = -22742388

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2699.7144
$ws.Range("I122").Value = 2825.5
$ws.Range("J122").Value = 1945
$ws.Range("K122").Value = 8476.5
$ws.Range("L122").Value = 5835
$ws.Range("M122").Value = -6026.5
$ws.Range("N122").Value = -10735

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1490.7693
$ws.Range("I61").Value = 922.5
$ws.Range("K61").Value = 922.5
$ws.Range("M61").Value = -720.5
$ws.Range("H68").Value = 2138
$ws.Range("I68").Value = 1996.6666
$ws.Range("J68").Value = 2350
$ws.Range("K68").Value = 1996.6666
$ws.Range("L68").Value = 2350
$ws.Range("M68").Value = -1247.6666
$ws.Range("N68").Value = -3848
$ws.Range("H71").Value = 2138
$ws.Range("I71").Value = 1996.6666
$ws.Range("J71").Value = 2350
$ws.Range("K71").Value = 9983.333000000001
$ws.Range("L71").Value = 11750
$ws.Range("M71").Value = -6239.333000000001
$ws.Range("N71").Value = -19238
$ws.Range("H74").Value = 146049.25
$ws.Range("I74").Value = 265098.5
$ws.Range("J74").Value = 27000
$ws.Range("K74").Value = 265098.5
$ws.Range("L74").Value = 27000
$ws.Range("M74").Value = -264100.5
$ws.Range("N74").Value = -28996
$ws.Range("H77").Value = 146049.25
$ws.Range("I77").Value = 265098.5
$ws.Range("J77").Value = 27000
$ws.Range("K77").Value = 795295.5
$ws.Range("L77").Value = 81000
$ws.Range("M77").Value = -790303.5
$ws.Range("N77").Value = -90984
$ws.Range("H113").Value = 1490.7693
$ws.Range("I113").Value = 922.5
$ws.Range("K113").Value = 922.5
$ws.Range("M113").Value = 1247.5
$ws.Range("H122").Value = 4064.6667
$ws.Range("J122").Value = 4064.6667
$ws.Range("L122").Value = 12194.0001
$ws.Range("N122").Value = -17094.0001
$ws.Range("H132").Value = 3957.8
$ws.Range("I132").Value = 2693.4783
$ws.Range("J132").Value = 8112
$ws.Range("K132").Value = 8080.4349
$ws.Range("L132").Value = 24336
$ws.Range("M132").Value = -5550.4349
$ws.Range("N132").Value = -29396

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 888.44446
$ws.Range("I96").Value = 888.44446
$ws.Range("K96").Value = 888.44446
$ws.Range("M96").Value = 484.55554
$ws.Range("H122").Value = 1133.3334
$ws.Range("I122").Value = 1133.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3400.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -950.0001999999999
$ws.Range("N122").ClearContents()
